# Update cryptocurrency price/volume snapshot (cryptos.xlsx) with the
# latest values pulled by the scheduled GitHub Actions job.
#
# Each row (A=rank idx, B=coin name, C=coinranking URL, D=price, E=1h
# volume %) is refreshed with newly scraped figures. A few coins also
# changed rank between runs (rows 39-41 reshuffled among Stacks / Hedera
# / Filecoin, and row 51 swapped from Cosmos to SuiNetwork), so those
# rows get their name/link updated too, not just the numbers.
#
# All of these sheet values are stored as literal text (inline/shared
# strings) in the workbook -- prices like "1.00" or "7.50" must stay
# text so the trailing zero survives; Excel's COM layer auto-coerces a
# plain `.Value = "7.50"` assignment on a numeric-looking string into the
# number 7.5 (dropping the formatting digit and flipping the cell to a
# numeric type). To avoid that, numeric-looking replacement values are
# written with a temporary "@" (Text) number format, then the format is
# cleared again with ClearFormats() so the cell keeps its original
# (unstyled) look once the literal text is safely stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '71.506.69' },
    @{ Cell = 'E2'; Value = '  +0.71%  ' },
    @{ Cell = 'D3'; Value = '3.816.89' },
    @{ Cell = 'E3'; Value = '  -0.89%  ' },
    @{ Cell = 'E4'; Value = '  -0.05%  ' },
    @{ Cell = 'D5'; Value = '702.19' },
    @{ Cell = 'E5'; Value = '  -0.74%  ' },
    @{ Cell = 'D6'; Value = '171.59' },
    @{ Cell = 'E6'; Value = '  -0.67%  ' },
    @{ Cell = 'D7'; Value = '3.815.05' },
    @{ Cell = 'E7'; Value = '  -0.88%  ' },
    @{ Cell = 'D8'; Value = '1.00' },
    @{ Cell = 'E8'; Value = '  +0.20%  ' },
    @{ Cell = 'D9'; Value = '0.527' },
    @{ Cell = 'E9'; Value = '  +0.08%  ' },
    @{ Cell = 'D10'; Value = '0.161' },
    @{ Cell = 'E10'; Value = '  -1.64%  ' },
    @{ Cell = 'D11'; Value = '7.50' },
    @{ Cell = 'E11'; Value = '  +2.19%  ' },
    @{ Cell = 'D12'; Value = '0.484' },
    @{ Cell = 'E12'; Value = '  +5.57%  ' },
    @{ Cell = 'D13'; Value = '0.0000252' },
    @{ Cell = 'E13'; Value = '  -1.87%  ' },
    @{ Cell = 'D14'; Value = '36.13' },
    @{ Cell = 'E14'; Value = '  -1.58%  ' },
    @{ Cell = 'D15'; Value = '4.454.37' },
    @{ Cell = 'E15'; Value = '  -1.00%  ' },
    @{ Cell = 'D16'; Value = '3.814.43' },
    @{ Cell = 'E16'; Value = '  -0.67%  ' },
    @{ Cell = 'D17'; Value = '71.446.33' },
    @{ Cell = 'E17'; Value = '  +0.60%  ' },
    @{ Cell = 'D18'; Value = '7.22' },
    @{ Cell = 'E18'; Value = '  +0.19%  ' },
    @{ Cell = 'D19'; Value = '17.59' },
    @{ Cell = 'E19'; Value = '  +1.29%  ' },
    @{ Cell = 'D20'; Value = '0.115' },
    @{ Cell = 'E20'; Value = '  -0.30%  ' },
    @{ Cell = 'D21'; Value = '515.55' },
    @{ Cell = 'E21'; Value = '  +3.56%  ' },
    @{ Cell = 'E22'; Value = '  -1.71%  ' },
    @{ Cell = 'D23'; Value = '0.717' },
    @{ Cell = 'E23'; Value = '  +0.00%  ' },
    @{ Cell = 'D24'; Value = '84.24' },
    @{ Cell = 'E24'; Value = '  -1.48%  ' },
    @{ Cell = 'D25'; Value = '0.0000143' },
    @{ Cell = 'E25'; Value = '  -2.90%  ' },
    @{ Cell = 'D26'; Value = '12.75' },
    @{ Cell = 'E26'; Value = '  +4.53%  ' },
    @{ Cell = 'D27'; Value = '3.959.21' },
    @{ Cell = 'E27'; Value = '  -1.04%  ' },
    @{ Cell = 'D28'; Value = '10.39' },
    @{ Cell = 'E28'; Value = '  -2.66%  ' },
    @{ Cell = 'E29'; Value = '  +0.13%  ' },
    @{ Cell = 'D30'; Value = '2.03' },
    @{ Cell = 'E30'; Value = '  -3.73%  ' },
    @{ Cell = 'D31'; Value = '3.02' },
    @{ Cell = 'E31'; Value = '  -5.84%  ' },
    @{ Cell = 'E32'; Value = '  -1.26%  ' },
    @{ Cell = 'D33'; Value = '7.36' },
    @{ Cell = 'E33'; Value = '  -2.23%  ' },
    @{ Cell = 'D34'; Value = '29.24' },
    @{ Cell = 'E34'; Value = '  -0.98%  ' },
    @{ Cell = 'E35'; Value = '  -3.73%  ' },
    @{ Cell = 'D36'; Value = '9.26' },
    @{ Cell = 'E36'; Value = '  +0.77%  ' },
    @{ Cell = 'D37'; Value = '3.775.52' },
    @{ Cell = 'E37'; Value = '  -0.83%  ' },
    @{ Cell = 'E38'; Value = '  -0.11%  ' },
    @{ Cell = 'B39'; Value = 'Stacks' },
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx' },
    @{ Cell = 'D39'; Value = '2.47' },
    @{ Cell = 'E39'; Value = '  +4.23%  ' },
    @{ Cell = 'B40'; Value = 'Hedera' },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' },
    @{ Cell = 'D40'; Value = '0.101' },
    @{ Cell = 'E40'; Value = '  -1.94%  ' },
    @{ Cell = 'B41'; Value = 'Filecoin' },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' },
    @{ Cell = 'D41'; Value = '6.34' },
    @{ Cell = 'E41'; Value = '  +4.75%  ' },
    @{ Cell = 'E42'; Value = '  -1.52%  ' },
    @{ Cell = 'D43'; Value = '3.27' },
    @{ Cell = 'E43'; Value = '  -2.23%  ' },
    @{ Cell = 'D45'; Value = '172.79' },
    @{ Cell = 'E45'; Value = '  +5.68%  ' },
    @{ Cell = 'E46'; Value = '  -0.09%  ' },
    @{ Cell = 'D47'; Value = '49.91' },
    @{ Cell = 'E48'; Value = '  -4.42%  ' },
    @{ Cell = 'D49'; Value = '427.69' },
    @{ Cell = 'E49'; Value = '  +2.69%  ' },
    @{ Cell = 'E50'; Value = '  -0.88%  ' },
    @{ Cell = 'B51'; Value = 'SuiNetwork' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui' },
    @{ Cell = 'D51'; Value = '1.17' },
    @{ Cell = 'E51'; Value = '  +4.23%  ' }
)

function Set-TextValue($Cell, $NewText) {
    # Plain decimals (optionally signed) are the only shapes Excel's COM
    # layer silently re-types as a number; thousands-dotted prices like
    # "71.506.69" and the "  +n.nn%  " volume strings never match this and
    # round-trip through .Value untouched.
    $looksNumeric = $NewText -match '^\s*[-+]?\d+(\.\d+)?\s*$'

    $range = $ws.Range($Cell)
    if ($looksNumeric) {
        # Force text storage so "1.00" / "7.50" don't get normalised to 1 / 7.5.
        $range.NumberFormat = "@"
        $range.Value = $NewText
        $range.ClearFormats()
    } else {
        $range.Value = $NewText
    }
}

foreach ($u in $updates) {
    Set-TextValue $u.Cell $u.Value
}
